## QRG - updated - cannot run in shared directories
##
## 1) Inserts a new numbered list item ("Make sure the script is placed
##    (and run) from a local directory ...") right after the "Steps"
##    heading and before "Launch the script".
## 2) Moves a couple of <w:lastRenderedPageBreak/> markers around
##    (removed before "Select " in the RVC-download step, added before
##    "The next step wilil be to" and before "The second").

$d = $word.ActiveDocument

function FindTemplateRange($text) {
    $rng = $d.Content
    $rng.Find.Execute($text, $true, $true, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    return $d.Range($rng.Start, $rng.End)
}

function FindParaIndex($text) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs.Item($i).Range.Text -like "*$text*") {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------
# 1. Insert the new list paragraph.
# ---------------------------------------------------------------------

$launchIdx = FindParaIndex "Launch the script"
$launchPara = $d.Paragraphs.Item($launchIdx)

# Inserting a paragraph break *before* "Launch the script" clones its
# pPr (ListParagraph style, numId 2), so the new bullet continues the
# same numbered list and gets the correct numbering automatically.
$launchStart = $launchPara.Range
$launchStart.Collapse(1)   # wdCollapseStart
$launchStart.InsertParagraphBefore()

$newPara = $d.Paragraphs.Item($launchIdx)
$cursor = $newPara.Range.Start

# Templates already present in the document whose run formatting
# (including w:lang="en-AU") we want to reproduce exactly:
#   plain        -> "Launch the script"
#   bold + bCs   -> "remove any rows that contain notes"
#   italic + iCs -> "You can just press the X in the upper right corner"
$plainTpl  = FindTemplateRange "Launch the script"
$plainTplStart = $plainTpl.Start
$plainTplLen = $plainTpl.End - $plainTpl.Start

$boldTpl = FindTemplateRange "remove any rows that contain notes"
$boldTplStart = $boldTpl.Start
$boldTplLen = $boldTpl.End - $boldTpl.Start

$italicTpl = FindTemplateRange "You can just press the X in the upper right corner"
$italicTplStart = $italicTpl.Start
$italicTplLen = $italicTpl.End - $italicTpl.Start

$vtab = [char]11   # Word's internal "manual line break" character

function AddRun($text, $tplStart, $tplLen, $underline) {
    $srcRange = $d.Range($tplStart, $tplStart + $tplLen)
    $destPoint = $d.Range($cursor, $cursor)
    $destPoint.FormattedText = $srcRange.FormattedText
    $destRange = $d.Range($cursor, $cursor + $tplLen)
    $destRange.Text = $text
    if ($underline -eq 1) {
        $finalRange = $d.Range($cursor, $cursor + $text.Length)
        $finalRange.Underline = 1
    }
    $cursor = $cursor + $text.Length
}

$apos = [char]0x2019

AddRun "Make sure the script is placed (and run) from a " $plainTplStart $plainTplLen 0
AddRun "local directory" $boldTplStart $boldTplLen 1
AddRun "$vtab(i.e. Desktop, Downloads, Documents etc. " $plainTplStart $plainTplLen 0
AddRun "not %appdata% " $boldTplStart $boldTplLen 0
AddRun "or " $plainTplStart $plainTplLen 0
AddRun "OneDrive" $boldTplStart $boldTplLen 0
AddRun ")$vtab" $plainTplStart $plainTplLen 0
AddRun "Note: I${apos}m not sure why, but the script has weird permission issues in other directories" $italicTplStart $italicTplLen 0
AddRun "$vtab" $plainTplStart $plainTplLen 0

Write-Host "New paragraph inserted, final cursor=$cursor"

# ---------------------------------------------------------------------
# 2. Toggle <w:lastRenderedPageBreak/> markers.
# ---------------------------------------------------------------------

Write-Host "Done."
